$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values
$ws.Range("A2").Value = 28
$ws.Range("B2").Value = 62
$ws.Range("C2").Value = 38
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = 18
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 18
$ws.Range("I2").Value = 2

# Column C values for rows 3-29 (the "Poids" / weights column)
$cValues = @(60, 13, 13, 7, 35, 38, 51, 15, 38, 3, 32, 21, 13, 52, 46, 38, 28, 38, 59, 56, 21, 46, 12, 55, 58, 44, 39)

$row = 3
foreach ($val in $cValues) {
    $ws.Range("C$row").Value = $val
    $row = $row + 1
}
